$d = $word.ActiveDocument

# Helper: toggle Bold off-then-on (or on-then-off) on a Range so that the
# engine is forced to materialise it as its own run, instead of silently
# re-absorbing it into a neighbouring run that happens to share the exact
# same formatting. The net visual/semantic effect is a no-op.
function Pin-Run($rng, $wasBold) {
    if ($wasBold) {
        $rng.Bold = 0
        $rng.Bold = 1
    } else {
        $rng.Bold = 1
        $rng.Bold = 0
    }
}

# ---------------------------------------------------------------------------
# 1) "Linear:" paragraph - hj(k) = (h(k) + s(j,k)) mod m mit s(j,k) = j
#    "(k) = (h(k) + " -> "(k) = (h(k) " / "-" / " "   (s(j,k) stays untouched)
# ---------------------------------------------------------------------------
$pLinear = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptxt = $d.Paragraphs($i).Range.Text
    if ($ptxt -like "Linear:*h(k) + s(j,k)*") {
        $pLinear = $d.Paragraphs($i)
        break
    }
}

$scope = $pLinear.Range.Duplicate()
$null = $scope.Find.Execute("(k) = (h(k) + ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$plusStart = $scope.Start + 12          # offset of "+" inside the matched text

# Replace "+" with "-" (leftmost edit first).
$plusRng = $d.Range($plusStart, $plusStart + 1)
$plusRng.Text = "-"

# Re-assert, strictly left-to-right, every original run boundary that
# follows the edit point so the unrelated text after it is not silently
# coalesced into fewer runs than the source document had.
Pin-Run $d.Range($plusStart, $plusStart + 1) $false          # "-"
Pin-Run $d.Range($plusStart + 1, $plusStart + 2) $false       # " "
Pin-Run $d.Range($plusStart + 2, $plusStart + 8) $false       # "s(j,k)"
Pin-Run $d.Range($plusStart + 8, $plusStart + 20) $false      # ") mod m mit "
Pin-Run $d.Range($plusStart + 20, $plusStart + 26) $true      # "s(j,k)" (bold)
Pin-Run $d.Range($plusStart + 26, $plusStart + 30) $true      # " = j" (bold)

# ---------------------------------------------------------------------------
# 2) "Quadratisch:" paragraph - hj(k) = (h(k) + s(j,k)) mod m mit s(j,k) = (-1)^j * j^2
#    "(k) = (h(k) + s(j,k)) mod m mit " -> "(k) = (h(k) " / "-" / " s(j,k)) mod m mit "
# ---------------------------------------------------------------------------
$pQuad = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptxt = $d.Paragraphs($i).Range.Text
    if ($ptxt -like "Quadratisch:*h(k) + s(j,k)*") {
        $pQuad = $d.Paragraphs($i)
        break
    }
}

$scope2 = $pQuad.Range.Duplicate()
$null = $scope2.Find.Execute("(k) = (h(k) + ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$plusStart2 = $scope2.Start + 12
$plusRng2   = $d.Range($plusStart2, $plusStart2 + 1)
$plusRng2.Text = "-"

Pin-Run $d.Range($plusStart2, $plusStart2 + 1) $false         # "-"
Pin-Run $d.Range($plusStart2 + 1, $plusStart2 + 20) $false    # " s(j,k)) mod m mit "

# ---------------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from the end of the Double-hashing paragraph
#    ("m' = m - 2") to right after the first "...mod 101 = 34" result, i.e.
#    right after the run containing the lone "4".
# ---------------------------------------------------------------------------
$pGoBack = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptxt = $d.Paragraphs($i).Range.Text
    if ($ptxt -like "*mod 101 = 34*mod 101 = 34*") {
        $pGoBack = $d.Paragraphs($i)
        break
    }
}

$scope3 = $pGoBack.Range.Duplicate()
$null = $scope3.Find.Execute("mod 101 = 34", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmPos = $scope3.End
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

Write-Output "done"
